# Generate Report for Handback
# Updates the localization-status workbook after a handback event:
#  - Status text moves from "Ready for handoff" to "Handed back: in sync with en-US"
#  - The zh-cn / de-de detail sheets get their "Latest Target File" and
#    "Latest Handback File" columns populated (with a hyperlink on the target file)
#    and a real "Latest Handback DateTime" stamp.
#  - A couple of columns are widened so the new/longer text is readable.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"
$mdFileName = "a923ed86-eea9-4a41-8447-fbda1f50c501.md"
$mdFileUrl  = "https://github.com/OpenLocalizationTestOrg/oltest/blob/92751b7b988795a79de0350678787246d8e98e05/e2e/a923ed86-eea9-4a41-8447-fbda1f50c501.md"

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (Overview!E2/F2 show per-language status, zh-cn!C2 / de-de!C2 show the
#    same status on each language's own sheet.)
# ---------------------------------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn sheet row 2: Latest Target File / Latest Handback File / Latest
#    Handback DateTime.
# ---------------------------------------------------------------------------
$wsZhCn.Range("I2").Value = $mdFileName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdFileUrl, "", "", $mdFileName) | Out-Null
$wsZhCn.Range("I2").Font.Color = 15570276

$wsZhCn.Range("J2").Value = "a923ed86-eea9-4a41-8447-fbda1f50c501.e5862f00626cb73e5e628373a2c44d53a29d366b.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-13 13:14:35"

# ---------------------------------------------------------------------------
# 3. de-de sheet row 2: Latest Target File / Latest Handback File / Latest
#    Handback DateTime.
# ---------------------------------------------------------------------------
$wsDeDe.Range("I2").Value = $mdFileName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdFileUrl, "", "", $mdFileName) | Out-Null
$wsDeDe.Range("I2").Font.Color = 15570276

$wsDeDe.Range("J2").Value = "a923ed86-eea9-4a41-8447-fbda1f50c501.e5862f00626cb73e5e628373a2c44d53a29d366b.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-13 13:14:45"

# ---------------------------------------------------------------------------
# 4. Column widths: widen the "Status" columns on Overview (E,F) and the
#    "Status" column on the zh-cn / de-de sheets (C) to fit the longer
#    status text, and widen "Latest Target File" / "Latest Handback File"
#    (I,J) on zh-cn / de-de to the standard wide (40) filename-column width.
# ---------------------------------------------------------------------------
$wideStatusWidth = 29.15   # lands on the same rendered width as the source workbook
$wideFileWidth   = 39.1666666666667  # renders as an exact width of 40

$wsOverview.Columns.Item(5).ColumnWidth = $wideStatusWidth   # E
$wsOverview.Columns.Item(6).ColumnWidth = $wideStatusWidth   # F

$wsZhCn.Columns.Item(3).ColumnWidth = $wideStatusWidth       # C
$wsZhCn.Columns.Item(9).ColumnWidth = $wideFileWidth         # I
$wsZhCn.Columns.Item(10).ColumnWidth = $wideFileWidth        # J

$wsDeDe.Columns.Item(3).ColumnWidth = $wideStatusWidth       # C
$wsDeDe.Columns.Item(9).ColumnWidth = $wideFileWidth         # I
$wsDeDe.Columns.Item(10).ColumnWidth = $wideFileWidth        # J

Write-Host "Handback report generated."
